$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the page count for "The Passionate Programmer" by Chad Fowler (row 13)
$ws.Range("C13").Value = 170

# Move the active selection to A10
$ws.Range("A10").Select()

$wb.Save()
